$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "booked" note for Sammy (row 5) and Harvey (row 8)
$ws.Range("C5").Value = "booked"
$ws.Range("C8").Value = "booked"

# Add a new row for James P
$ws.Range("B13").Value = "James P"
$ws.Range("C13").Value = "booked"

# Update the active selection to match the end-user's final selection
$ws.Range("F10").Select()
